# Auto-generated edit script: updates market-price derived values in the
# Balmung_Profits workbook per the scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

### Sheet: ALC (42 cell updates)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 3613.5715
$ws.Range("I4").Value = 1968.4615
$ws.Range("K4").Value = 1968.4615
$ws.Range("M4").Value = -1854.4615
$ws.Range("H9").Value = 8342093.5
$ws.Range("J9").Value = 28046.334
$ws.Range("L9").Value = 28046.334
$ws.Range("N9").Value = -28384.334
$ws.Range("H18").Value = 5232.8335
$ws.Range("I18").Value = 5232.8335
$ws.Range("K18").Value = 5232.8335
$ws.Range("M18").Value = -4948.8335
$ws.Range("H33").Value = 3164588.8
$ws.Range("I33").Value = 3917435
$ws.Range("K33").Value = 3917435
$ws.Range("M33").Value = -3917206
$ws.Range("H41").Value = 1568.375
$ws.Range("I41").Value = 1139.6364
$ws.Range("J41").Value = 1931.1538
$ws.Range("K41").Value = 1139.6364
$ws.Range("L41").Value = 1931.1538
$ws.Range("M41").Value = -699.6364000000001
$ws.Range("N41").Value = -2811.1538
$ws.Range("H98").Value = 6053.2856
$ws.Range("I98").Value = 4666
$ws.Range("K98").Value = 4666
$ws.Range("M98").Value = -3168
$ws.Range("H106").Value = 55558620
$ws.Range("I106").Value = 58826332
$ws.Range("K106").Value = 58826332
$ws.Range("M106").Value = -58825701
$ws.Range("H122").Value = 6053.2856
$ws.Range("I122").Value = 4666
$ws.Range("K122").Value = 13998
$ws.Range("M122").Value = -11548
$ws.Range("H135").Value = 1151.2646
$ws.Range("I135").Value = 801.12
$ws.Range("J135").Value = 2123.889
$ws.Range("K135").Value = 7210.08
$ws.Range("L135").Value = 19115.001
$ws.Range("M135").Value = -4675.08
$ws.Range("N135").Value = -24185.001

### Sheet: ARM (22 cell updates)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 324345.56
$ws.Range("I74").Value = 1222.6981
$ws.Range("J74").Value = 1394690
$ws.Range("K74").Value = 1222.6981
$ws.Range("L74").Value = 1394690
$ws.Range("M74").Value = -348.6981000000001
$ws.Range("N74").Value = -1396438
$ws.Range("H77").Value = 324345.56
$ws.Range("I77").Value = 1222.6981
$ws.Range("J77").Value = 1394690
$ws.Range("K77").Value = 6113.4905
$ws.Range("L77").Value = 6973450
$ws.Range("M77").Value = -1745.4905
$ws.Range("N77").Value = -6982186
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = $null

### Sheet: BSM (64 cell updates)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1695.7142
$ws.Range("I22").Value = 1194
$ws.Range("J22").Value = 2950
$ws.Range("K22").Value = 1194
$ws.Range("L22").Value = 2950
$ws.Range("M22").Value = -1021
$ws.Range("N22").Value = -3296
$ws.Range("H36").Value = 1425.2858
$ws.Range("I36").Value = 1496.1666
$ws.Range("J36").Value = 1000
$ws.Range("K36").Value = 1496.1666
$ws.Range("L36").Value = 1000
$ws.Range("M36").Value = -962.1666
$ws.Range("N36").Value = -2068
$ws.Range("H68").Value = 60000
$ws.Range("J68").Value = 60000
$ws.Range("L68").Value = 60000
$ws.Range("N68").Value = -61622
$ws.Range("H71").Value = 60000
$ws.Range("J71").Value = 60000
$ws.Range("L71").Value = 180000
$ws.Range("N71").Value = -188112
$ws.Range("H86").Value = 3157.2856
$ws.Range("I86").Value = 1133.8889
$ws.Range("K86").Value = 1133.8889
$ws.Range("M86").Value = -10.88889999999992
$ws.Range("H89").Value = 3157.2856
$ws.Range("I89").Value = 1133.8889
$ws.Range("K89").Value = 5669.4445
$ws.Range("M89").Value = -53.44449999999961
$ws.Range("H94").Value = 14142.571
$ws.Range("I94").Value = 15598.8
$ws.Range("K94").Value = 15598.8
$ws.Range("M94").Value = -15147.8
$ws.Range("H95").Value = 60000
$ws.Range("J95").Value = 60000
$ws.Range("L95").Value = 60000
$ws.Range("N95").Value = -65492
$ws.Range("H97").Value = 25013.25
$ws.Range("I97").Value = 4276.5
$ws.Range("J97").Value = 45750
$ws.Range("K97").Value = 4276.5
$ws.Range("L97").Value = 45750
$ws.Range("M97").Value = -3285.5
$ws.Range("N97").Value = -47732
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = $null
$ws.Range("H103").Value = 13519.2
$ws.Range("J103").Value = 13519.2
$ws.Range("L103").Value = 13519.2
$ws.Range("N103").Value = -15863.2
$ws.Range("H107").Value = 8612.763000000001
$ws.Range("I107").Value = 9968.666999999999
$ws.Range("K107").Value = 9968.666999999999
$ws.Range("M107").Value = -8048.666999999999
$ws.Range("H134").Value = 34618020
$ws.Range("I134").Value = 2698.4443
$ws.Range("J134").Value = 112502500
$ws.Range("K134").Value = 8095.3329
$ws.Range("L134").Value = 337507500
$ws.Range("M134").Value = -5560.3329
$ws.Range("N134").Value = -337512570

### Sheet: CRP (46 cell updates)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2741.2222
$ws.Range("I31").Value = 2507.5293
$ws.Range("J31").Value = 2950.3157
$ws.Range("K31").Value = 2507.5293
$ws.Range("L31").Value = 2950.3157
$ws.Range("M31").Value = -2212.5293
$ws.Range("N31").Value = -3540.3157
$ws.Range("H34").Value = 2741.2222
$ws.Range("I34").Value = 2507.5293
$ws.Range("J34").Value = 2950.3157
$ws.Range("K34").Value = 2507.5293
$ws.Range("L34").Value = 2950.3157
$ws.Range("M34").Value = -2305.5293
$ws.Range("N34").Value = -3354.3157
$ws.Range("H86").Value = 11395.05
$ws.Range("I86").Value = 7194.385
$ws.Range("J86").Value = 19196.285
$ws.Range("K86").Value = 7194.385
$ws.Range("L86").Value = 19196.285
$ws.Range("M86").Value = -6071.385
$ws.Range("N86").Value = -21442.285
$ws.Range("H89").Value = 11395.05
$ws.Range("I89").Value = 7194.385
$ws.Range("J89").Value = 19196.285
$ws.Range("K89").Value = 35971.925
$ws.Range("L89").Value = 95981.425
$ws.Range("M89").Value = -30355.925
$ws.Range("N89").Value = -107213.425
$ws.Range("H132").Value = 24387.238
$ws.Range("J132").Value = 3331.3125
$ws.Range("L132").Value = 9993.9375
$ws.Range("N132").Value = -15053.9375
$ws.Range("H134").Value = 1235.6731
$ws.Range("I134").Value = 1206.2632
$ws.Range("J134").Value = 1315.5
$ws.Range("K134").Value = 3618.7896
$ws.Range("L134").Value = 3946.5
$ws.Range("M134").Value = -1083.7896
$ws.Range("N134").Value = -9016.5
$ws.Range("H141").Value = 486788.72
$ws.Range("I141").Value = 373325.66
$ws.Range("J141").Value = 529337.4
$ws.Range("K141").Value = 373325.66
$ws.Range("L141").Value = 529337.4
$ws.Range("M141").Value = -368145.66
$ws.Range("N141").Value = -539697.4

### Sheet: CUL (33 cell updates)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = $null
$ws.Range("N64").Value = $null
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = $null
$ws.Range("N67").Value = $null
$ws.Range("H81").Value = 23817796
$ws.Range("J81").Value = 10264
$ws.Range("L81").Value = 30792
$ws.Range("N81").Value = -33038
$ws.Range("H84").Value = 23817796
$ws.Range("J84").Value = 10264
$ws.Range("L84").Value = 92376
$ws.Range("N84").Value = -103608
$ws.Range("H97").Value = 450.7143
$ws.Range("I97").Value = 407.5
$ws.Range("J97").Value = 710
$ws.Range("K97").Value = 1222.5
$ws.Range("L97").Value = 2130
$ws.Range("M97").Value = -726.5
$ws.Range("N97").Value = -3122
$ws.Range("H140").Value = 12502101
$ws.Range("I140").Value = 13159580
$ws.Range("K140").Value = 39478740
$ws.Range("M140").Value = -39473560

### Sheet: GSM (7 cell updates)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2904.8572
$ws.Range("I126").Value = 2889
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8667
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -6197
$ws.Range("N126").Value = -13940

### Sheet: LTW (30 cell updates)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 23703.5
$ws.Range("J2").Value = 61110.5
$ws.Range("L2").Value = 61110.5
$ws.Range("N2").Value = -61334.5
$ws.Range("H16").Value = 2230.0715
$ws.Range("I16").Value = 2093.9614
$ws.Range("K16").Value = 2093.9614
$ws.Range("M16").Value = -1923.9614
$ws.Range("H40").Value = 3303.6296
$ws.Range("I40").Value = 2249.389
$ws.Range("J40").Value = 5412.1113
$ws.Range("K40").Value = 2249.389
$ws.Range("L40").Value = 5412.1113
$ws.Range("M40").Value = -2113.389
$ws.Range("N40").Value = -5684.1113
$ws.Range("H100").Value = 3130.75
$ws.Range("I100").Value = 2355.6
$ws.Range("J100").Value = 4422.6665
$ws.Range("K100").Value = 2355.6
$ws.Range("L100").Value = 4422.6665
$ws.Range("M100").Value = -1814.6
$ws.Range("N100").Value = -5504.6665
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = $null
$ws.Range("H132").Value = 3089.0625
$ws.Range("I132").Value = 2922.037
$ws.Range("K132").Value = 8766.110999999999
$ws.Range("M132").Value = -6236.110999999999

### Sheet: WVR (15 cell updates)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 33567.65
$ws.Range("J14").Value = 38443.46
$ws.Range("L14").Value = 38443.46
$ws.Range("N14").Value = -38779.46
$ws.Range("H132").Value = 1831.0883
$ws.Range("I132").Value = 1630.4482
$ws.Range("J132").Value = 2994.8
$ws.Range("K132").Value = 4891.3446
$ws.Range("L132").Value = 8984.400000000001
$ws.Range("M132").Value = -2361.3446
$ws.Range("N132").Value = -14044.4
$ws.Range("H136").Value = 18357.75
$ws.Range("I136").Value = 25185.021
$ws.Range("K136").Value = 75555.06299999999
$ws.Range("M136").Value = -73005.06299999999

